$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Hora entrada" / "Hora salida" header cells and the per-row time
# values are being unlinked/cleared (traslados y desvinculacion) while
# keeping their existing formatting/styles intact.
$ws.Range("B1:C1").ClearContents()
$ws.Range("B2:C2").ClearContents()
$ws.Range("B3:C3").ClearContents()

# Reflect the resulting selection left behind by the edit.
$ws.Range("B1:C5").Select()
